$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text, possibly numeric-looking values (e.g. "570.75",
# "5.10") that must stay exactly as scraped text (preserving trailing zeros,
# thousands-dot grouping, subscript digits, etc.). Force each target cell to
# Text format before assigning so Excel does not silently convert it to a
# floating point number, then restore the default "Normal" style so no stray
# formatting remains on the cell.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

$ws.Range('D2').Value = '61.035.63'
$ws.Range('D3').Value = '2.422.74'
$ws.Range('D5').Value = '570.75'
$ws.Range('D6').Value = '139.91'
$ws.Range('D8').Value = '0.526'
$ws.Range('D9').Value = '2.408.43'
$ws.Range('D12').Value = '5.10'
$ws.Range('D14').Value = '26.19'
$ws.Range('D17').Value = '60.885.74'
$ws.Range('D18').Value = '2.418.72'
$ws.Range('D19').Value = '7.60'
$ws.Range('D21').Value = '323.44'
$ws.Range('D22').Value = '4.06'
$ws.Range('D26').Value = '64.88'
$ws.Range('D27').Value = '588.43'
$ws.Range('D29').Value = '2.546.88'
$ws.Range('D30').Value = '0.0₃0936'
$ws.Range('D32').Value = '1.35'
$ws.Range('D38').Value = '151.94'
$ws.Range('D39').Value = '0.368'
$ws.Range('D40').Value = '18.31'
$ws.Range('D44').Value = '41.26'
$ws.Range('D45').Value = '2.36'
$ws.Range('D46').Value = '0.0₆0292'
$ws.Range('D47').Value = '142.24'
$ws.Range('D48').Value = '3.52'
$ws.Range('D50').Value = '19.60'

$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'

# Column E values are plain percentage text (e.g. "  -1.92%  ") and are safe to
# assign directly since Excel will not reinterpret them as numbers.
$ws.Range('E2').Value = '  -1.92%  '
$ws.Range('E3').Value = '  -1.01%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('E5').Value = '  -2.10%  '
$ws.Range('E6').Value = '  -2.24%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('E8').Value = '  -0.82%  '
$ws.Range('E9').Value = '  -1.37%  '
$ws.Range('E10').Value = '  -1.09%  '
$ws.Range('E11').Value = '  -0.19%  '
$ws.Range('E12').Value = '  -2.00%  '
$ws.Range('E13').Value = '  -1.84%  '
$ws.Range('E14').Value = '  -0.94%  '
$ws.Range('E15').Value = '  -3.95%  '
$ws.Range('E17').Value = '  -2.00%  '
$ws.Range('E18').Value = '  -0.90%  '
$ws.Range('E19').Value = '  +7.18%  '
$ws.Range('E20').Value = '  -1.17%  '
$ws.Range('E21').Value = '  -1.22%  '
$ws.Range('E22').Value = '  -1.27%  '
$ws.Range('E23').Value = '  +1.37%  '
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('E25').Value = '  -2.98%  '
$ws.Range('E26').Value = '  -1.22%  '
$ws.Range('E27').Value = '  -0.26%  '
$ws.Range('E28').Value = '  -9.03%  '
$ws.Range('E29').Value = '  -0.91%  '
$ws.Range('E30').Value = '  -3.32%  '
$ws.Range('E31').Value = '  -0.75%  '
$ws.Range('E32').Value = '  -4.27%  '
$ws.Range('E33').Value = '  -1.46%  '
$ws.Range('E34').Value = '  -1.96%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('E36').Value = '  -1.31%  '
$ws.Range('E37').Value = '  -5.52%  '
$ws.Range('E38').Value = '  -0.93%  '
$ws.Range('E39').Value = '  -2.33%  '
$ws.Range('E40').Value = '  -0.55%  '
$ws.Range('E41').Value = '  -2.19%  '
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('E43').Value = '  -1.20%  '
$ws.Range('E44').Value = '  -3.71%  '
$ws.Range('E45').Value = '  -5.09%  '
$ws.Range('E46').Value = '  +14.12%  '
$ws.Range('E47').Value = '  -0.19%  '
$ws.Range('E48').Value = '  -3.37%  '
$ws.Range('E49').Value = '  -2.32%  '
$ws.Range('E50').Value = '  -1.11%  '
$ws.Range('E51').Value = '  -3.45%  '
